$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header
$ws.Range("A1").Value = "Cluster Name"
$ws.Range("B1").Value = "Active Cases"

# Data rows (Cluster Name, Active Cases)
$data = @(
    @('3364 Assisi Centre Aged Care Rosanna', 29),
    @('3376 Royal Freemasons Coppin Centre Melbourne', 16),
    @('3633 Lifeview Emerald Glades Aged Care Emerald', 10),
    @('3825 TLC Forest Lodge Residential Aged Care Frankston North', 10),
    @('3961 Heritage Care Water Gardens Aged Care Facility Sydenham', 23),
    @('4167 Royal Freemasons Centennial Lodge Wantirna South', 22),
    @('4282 Villa Maria Catholic Homes (VMCH) Wantirna At-Home Aged Care Stud Road Wantirna', 11),
    @('AG Industries Pty Ltd Factory Thomastown', 17),
    @('Antonine College Cedar Campus Coburg', 12),
    @('Australian Meat Group Abattoir Dandenong South', 10),
    @('Bacchus Marsh Childcare and Kindergarten Centre Bacchus Marsh', 36),
    @('Baden Powell College Tarneit', 14),
    @('Chairo Christian School Pakenham Campus Pakenham', 10),
    @('Covenant College Bell Post Hill', 22),
    @('Devon Meadows Primary School Devon Meadows', 12),
    @('Domestic Freight Vessel Siem Aquamarine Geelong', 10),
    @('Gilly''s Early Learning Centre Balaclava', 10),
    @('Gladstone Park Secondary College 29 Oct Gladstone Park', 13),
    @('Guardian Childcare & Education Moorabbin', 11),
    @('Hamlyn Banks Primary School Hamlyn Heights', 10),
    @('Hazel Glen College Doreen', 11),
    @('Hazelwood North Primary School Hazelwood North', 28),
    @('Islamic College of Melbourne Tarneit Oct Nov', 27),
    @('Lyndhurst Primary School Lyndhurst', 11),
    @('Master Poultry Group West Footscray', 14),
    @('Minaret College Officer Campus Officer', 12),
    @('Morwell Park Primary School Morwell', 57),
    @('Nido Early School Woodend', 12),
    @('Northern Bay College Goldsworthy 9-12 Campus Corio', 10),
    @('Northern Bay College Wexford Campus Corio', 31),
    @('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 14),
    @('Northern Health The Northern Hospital Epping', 14),
    @('Oakleigh South Primary School Oakleigh South', 18),
    @('Our Lady''s Catholic Primary School Wangaratta', 10),
    @('Pentland Primary School', 10),
    @('Rutherglen Motor Inn and Walkabout Motel Rutherglen', 22),
    @('Sirius College Ibrahim Dellal Campus Sunshine', 11),
    @('Sirius College Shepparton Campus Shepparton', 15),
    @('Smartie Pants Early Learning and Development Diamond Creek', 20),
    @('Social Gathering Woodvale 30 Oct', 10),
    @('Society Restaurant Melbourne', 21),
    @('St Brendans Primary School Shepparton', 11),
    @('St Clare''s Primary School Officer', 10),
    @('St Georges Road Primary School Shepparton', 14),
    @('St Joseph''s School Quarry Hill', 32),
    @('St Josephs Catholic Primary School Warragul', 10),
    @('St Louis de Montfort''s School Aspendale', 13),
    @('St Paul''s Primary School Sunshine West', 13),
    @('St Vincents Hospital Emergency Department Melbourne', 10),
    @('Stevensville Primary School St Albans', 10),
    @('Stockdale Road Primary School Traralgon', 34),
    @('Story House Early Learning Epping October', 13),
    @('Sunbury Primary School Sunbury', 12),
    @('TUROSI PTY LTD Thomastown', 10),
    @('Templestowe Park Primary School Templestowe', 31),
    @('The Lake Primary School Cabarita', 19),
    @('The Royal Children''s Hospital Melbourne Emergency Department Parkville', 10),
    @('Top Yard Rooftop Melbourne', 13),
    @('Tucker Road Bentleigh Primary School Bentleigh', 10),
    @('Warragul Regional College Warragul', 18),
    @('Werribee Mercy Hospital Emergency Department', 22),
    @('Western Health Sunshine Hospital Emergency Department St Albans', 13),
    @('Wodonga Primary School Wodonga', 13),
    @('Wodonga South Primary School Wodonga', 10),
    @('Wyndham Christian College Wyndham Vale', 12),
    @('Yeshivah College St Kilda East', 19)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove any leftover rows beyond the new data range (in case sheet had more rows before)
$lastRow = $startRow + $data.Count - 1
$ws.Range("A" + ($lastRow + 1) + ":B10000").ClearContents()
